$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new header columns for p75 stats (AR, AS, AT)
$ws.Cells.Item(1, 44).Value = "p75CellsScored"
$ws.Cells.Item(1, 45).Value = "p75CellsScoredHigh"
$ws.Cells.Item(1, 46).Value = "p75CellsScoredLow"

# Per-team p75 values computed from match data (cellsScored / High / Low)
$ws.Cells.Item(2, 44).Value = 8.699999999999999
$ws.Cells.Item(2, 45).Value = 8.699999999999999
$ws.Cells.Item(2, 46).Value = 0
$ws.Cells.Item(3, 44).Value = 2.3
$ws.Cells.Item(3, 45).Value = 2.3
$ws.Cells.Item(3, 46).Value = 0
$ws.Cells.Item(4, 44).Value = 5.5
$ws.Cells.Item(4, 45).Value = 5.5
$ws.Cells.Item(4, 46).Value = 0
$ws.Cells.Item(5, 44).Value = 9
$ws.Cells.Item(5, 45).Value = 9
$ws.Cells.Item(5, 46).Value = 0
$ws.Cells.Item(6, 44).Value = 2.7
$ws.Cells.Item(6, 45).Value = 0
$ws.Cells.Item(6, 46).Value = 2.7
$ws.Cells.Item(7, 44).Value = 10.7
$ws.Cells.Item(7, 45).Value = 10.7
$ws.Cells.Item(7, 46).Value = 0
$ws.Cells.Item(8, 44).Value = 0
$ws.Cells.Item(8, 45).Value = 0
$ws.Cells.Item(8, 46).Value = 0
$ws.Cells.Item(9, 44).Value = 3
$ws.Cells.Item(9, 45).Value = 0
$ws.Cells.Item(9, 46).Value = 3
$ws.Cells.Item(10, 44).Value = 3
$ws.Cells.Item(10, 45).Value = 2
$ws.Cells.Item(10, 46).Value = 1
$ws.Cells.Item(11, 44).Value = 0
$ws.Cells.Item(11, 45).Value = 0
$ws.Cells.Item(11, 46).Value = 0
$ws.Cells.Item(12, 44).Value = 0
$ws.Cells.Item(12, 45).Value = 0
$ws.Cells.Item(12, 46).Value = 0
$ws.Cells.Item(13, 44).Value = 10.3
$ws.Cells.Item(13, 45).Value = 1.3
$ws.Cells.Item(13, 46).Value = 9.300000000000001
$ws.Cells.Item(14, 44).Value = 0
$ws.Cells.Item(14, 45).Value = 0
$ws.Cells.Item(14, 46).Value = 0
$ws.Cells.Item(15, 44).Value = 6.7
$ws.Cells.Item(15, 45).Value = 0
$ws.Cells.Item(15, 46).Value = 6.7
$ws.Cells.Item(16, 44).Value = 8.699999999999999
$ws.Cells.Item(16, 45).Value = 8.699999999999999
$ws.Cells.Item(16, 46).Value = 0
$ws.Cells.Item(17, 44).Value = 11.7
$ws.Cells.Item(17, 45).Value = 11.7
$ws.Cells.Item(17, 46).Value = 0
$ws.Cells.Item(18, 44).Value = 3.7
$ws.Cells.Item(18, 45).Value = 3.7
$ws.Cells.Item(18, 46).Value = 0
$ws.Cells.Item(19, 44).Value = 1
$ws.Cells.Item(19, 45).Value = 1
$ws.Cells.Item(19, 46).Value = 0
$ws.Cells.Item(20, 44).Value = 2.7
$ws.Cells.Item(20, 45).Value = 2.7
$ws.Cells.Item(20, 46).Value = 0
$ws.Cells.Item(21, 44).Value = 4
$ws.Cells.Item(21, 45).Value = 4
$ws.Cells.Item(21, 46).Value = 0
$ws.Cells.Item(22, 44).Value = 6.3
$ws.Cells.Item(22, 45).Value = 6.3
$ws.Cells.Item(22, 46).Value = 0
$ws.Cells.Item(23, 44).Value = 0
$ws.Cells.Item(23, 45).Value = 0
$ws.Cells.Item(23, 46).Value = 0
$ws.Cells.Item(24, 44).Value = 10
$ws.Cells.Item(24, 45).Value = 10
$ws.Cells.Item(24, 46).Value = 0
$ws.Cells.Item(25, 44).Value = 15.3
$ws.Cells.Item(25, 45).Value = 15.3
$ws.Cells.Item(25, 46).Value = 0
$ws.Cells.Item(26, 44).Value = 8.5
$ws.Cells.Item(26, 45).Value = 8.5
$ws.Cells.Item(26, 46).Value = 0
$ws.Cells.Item(27, 44).Value = 3
$ws.Cells.Item(27, 45).Value = 0
$ws.Cells.Item(27, 46).Value = 3
$ws.Cells.Item(28, 44).Value = 8.300000000000001
$ws.Cells.Item(28, 45).Value = 8.300000000000001
$ws.Cells.Item(28, 46).Value = 0
$ws.Cells.Item(29, 44).Value = 0
$ws.Cells.Item(29, 45).Value = 0
$ws.Cells.Item(29, 46).Value = 0
$ws.Cells.Item(30, 44).Value = 8
$ws.Cells.Item(30, 45).Value = 8
$ws.Cells.Item(30, 46).Value = 0
$ws.Cells.Item(31, 44).Value = 0.5
$ws.Cells.Item(31, 45).Value = 0.5
$ws.Cells.Item(31, 46).Value = 0
$ws.Cells.Item(32, 44).Value = 4
$ws.Cells.Item(32, 45).Value = 4
$ws.Cells.Item(32, 46).Value = 0
$ws.Cells.Item(33, 44).Value = 17
$ws.Cells.Item(33, 45).Value = 17
$ws.Cells.Item(33, 46).Value = 0
$ws.Cells.Item(34, 44).Value = 0
$ws.Cells.Item(34, 45).Value = 0
$ws.Cells.Item(34, 46).Value = 0
$ws.Cells.Item(35, 44).Value = 3.7
$ws.Cells.Item(35, 45).Value = 0
$ws.Cells.Item(35, 46).Value = 3.7
$ws.Cells.Item(36, 44).Value = 1.5
$ws.Cells.Item(36, 45).Value = 1.5
$ws.Cells.Item(36, 46).Value = 0
$ws.Cells.Item(37, 44).Value = 3
$ws.Cells.Item(37, 45).Value = 3
$ws.Cells.Item(37, 46).Value = 0.3
$ws.Cells.Item(38, 44).Value = 2
$ws.Cells.Item(38, 45).Value = 2
$ws.Cells.Item(38, 46).Value = 0
$ws.Cells.Item(39, 44).Value = 3.5
$ws.Cells.Item(39, 45).Value = 3.5
$ws.Cells.Item(39, 46).Value = 0
$ws.Cells.Item(40, 44).Value = 3.5
$ws.Cells.Item(40, 45).Value = 3
$ws.Cells.Item(40, 46).Value = 1
$ws.Cells.Item(41, 44).Value = 0.5
$ws.Cells.Item(41, 45).Value = 0
$ws.Cells.Item(41, 46).Value = 0.5
$ws.Cells.Item(42, 44).Value = 6
$ws.Cells.Item(42, 45).Value = 6
$ws.Cells.Item(42, 46).Value = 0
$ws.Cells.Item(43, 44).Value = 2.5
$ws.Cells.Item(43, 45).Value = 1
$ws.Cells.Item(43, 46).Value = 1.5
$ws.Cells.Item(44, 44).Value = 0.3
$ws.Cells.Item(44, 45).Value = 0.3
$ws.Cells.Item(44, 46).Value = 0
$ws.Cells.Item(45, 44).Value = 5.5
$ws.Cells.Item(45, 45).Value = 0
$ws.Cells.Item(45, 46).Value = 5.5
$ws.Cells.Item(46, 44).Value = 3.7
$ws.Cells.Item(46, 45).Value = 3.3
$ws.Cells.Item(46, 46).Value = 0.7
$ws.Cells.Item(47, 44).Value = 1.3
$ws.Cells.Item(47, 45).Value = 1.3
$ws.Cells.Item(47, 46).Value = 0
$ws.Cells.Item(48, 44).Value = 4
$ws.Cells.Item(48, 45).Value = 4
$ws.Cells.Item(48, 46).Value = 0
$ws.Cells.Item(49, 44).Value = 5.3
$ws.Cells.Item(49, 45).Value = 5.3
$ws.Cells.Item(49, 46).Value = 0
$ws.Cells.Item(50, 44).Value = 6.7
$ws.Cells.Item(50, 45).Value = 0.3
$ws.Cells.Item(50, 46).Value = 6.7
$ws.Cells.Item(51, 44).Value = 2.3
$ws.Cells.Item(51, 45).Value = 2.3
$ws.Cells.Item(51, 46).Value = 0.7
$ws.Cells.Item(52, 44).Value = 10
$ws.Cells.Item(52, 45).Value = 10
$ws.Cells.Item(52, 46).Value = 0
$ws.Cells.Item(53, 44).Value = 6.3
$ws.Cells.Item(53, 45).Value = 6.3
$ws.Cells.Item(53, 46).Value = 0
$ws.Cells.Item(54, 44).Value = 3
$ws.Cells.Item(54, 45).Value = 3
$ws.Cells.Item(54, 46).Value = 0
$ws.Cells.Item(55, 44).Value = 13.5
$ws.Cells.Item(55, 45).Value = 0
$ws.Cells.Item(55, 46).Value = 13.5
$ws.Cells.Item(56, 44).Value = 6.3
$ws.Cells.Item(56, 45).Value = 6.3
$ws.Cells.Item(56, 46).Value = 0
$ws.Cells.Item(57, 44).Value = 8
$ws.Cells.Item(57, 45).Value = 8
$ws.Cells.Item(57, 46).Value = 0
$ws.Cells.Item(58, 44).Value = 3.7
$ws.Cells.Item(58, 45).Value = 3.7
$ws.Cells.Item(58, 46).Value = 0
$ws.Cells.Item(59, 44).Value = 0
$ws.Cells.Item(59, 45).Value = 0
$ws.Cells.Item(59, 46).Value = 0
$ws.Cells.Item(60, 44).Value = 2.7
$ws.Cells.Item(60, 45).Value = 2.7
$ws.Cells.Item(60, 46).Value = 0
$ws.Cells.Item(61, 44).Value = 4.7
$ws.Cells.Item(61, 45).Value = 2.7
$ws.Cells.Item(61, 46).Value = 3
